$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.147.95"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.827.53"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value2 = "'242.09"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value2 = "'0.6217"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value2 = "'0.07366"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value2 = "'23.10"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D11").Value2 = "'0.07668"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.832.44"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value2 = "'4.944"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value2 = "'0.6648"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value2 = "'82.20"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value2 = "'0.000008981"
$ws.Range("E16").Value = "  -3.52%  "
$ws.Range("D17").Value2 = "'5.838"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "29.115.66"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "2.075.91"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value2 = "'7.332"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").Value2 = "'0.9990"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value2 = "'158.19"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value2 = "'0.1410"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value2 = "'8.490"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value2 = "'17.61"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value2 = "'0.06024"
$ws.Range("E29").Value = "  +7.47%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value2 = "'1.480"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value2 = "'4.087"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value2 = "'4.066"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").Value2 = "'1.205"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value2 = "'1.864"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value2 = "'0.7286"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value2 = "'1.139"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").Value2 = "'2.608"
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("D38").Value2 = "'2.838"
$ws.Range("D39").Value = "1.221.74"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value2 = "'0.01749"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value2 = "'6.275"
$ws.Range("E41").Value = "  -5.35%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value2 = "'0.9164"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value2 = "'101.71"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "1.981.29"
$ws.Range("D46").Value2 = "'64.82"
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value2 = "'0.5046"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value2 = "'0.4014"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value2 = "'9.088"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value2 = "'0.1131"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value2 = "'0.05749"
$ws.Range("E51").Value = "  -1.58%  "
